# Generate Report for Handoff
# Updates the "3de4c1a0-97db-4ece-afbd-afe4c8cfb490.md" row across the
# Overview / zh-cn / de-de sheets to reflect a fresh "Ready for handoff"
# status with updated timestamps and an error detail message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4e8a9593bed6f9d469097abe4f130773eecd4a4/e2e/3de4c1a0-97db-4ece-afbd-afe4c8cfb490.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/654d2648e151174ecf25b07a2cdc8d03a3f8c7a2/e2e/3de4c1a0-97db-4ece-afbd-afe4c8cfb490.md."

$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-03 02:52:40"
$wsOverview.Range("G3").NumberFormat = $dateTimeFormat

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-09-03 02:52:36"
$wsZhCn.Range("H3").NumberFormat = $dateTimeFormat
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-09-03 02:52:40"
$wsDeDe.Range("H3").NumberFormat = $dateTimeFormat
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
